# Update scraped schedule data (Línea 141) for the new scrape run at 04:32:18.
# Refreshes timestamps/counters on each sheet and rewrites/extends the data
# rows on "LP1912" (now 25 rows) and "LP1912-215" (now 5 rows); "6203-6173"
# only gets its "Última actualización" timestamp bumped.
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 04:32:18'
$ws.Range('A3').Value = 'Total filas: 25'
$ws.Range('A7').Value = '04:32:18'
$ws.Range('B7').Value = '04:33'
$ws.Range('C7').Value = '15_ABASTO'
$ws.Range('D7').Value = 1
$ws.Range('A8').Value = '04:13:31'
$ws.Range('C8').Value = '215_EL PELIGRO'
$ws.Range('D8').Value = 33
$ws.Range('A9').Value = '03:52:04'
$ws.Range('B9').Value = '04:46'
$ws.Range('C9').Value = '215A_EL PATO'
$ws.Range('D9').Value = 54
$ws.Range('A10').Value = '04:32:18'
$ws.Range('B10').Value = '04:47'
$ws.Range('C10').Value = '215_EL PELIGRO'
$ws.Range('D10').Value = 15
$ws.Range('A11').Value = '04:32:18'
$ws.Range('B11').Value = '04:53'
$ws.Range('C11').Value = '11_ETCHEVERRY'
$ws.Range('D11').Value = 21
$ws.Range('B12').Value = '05:11'
$ws.Range('C12').Value = '17_ROMERO'
$ws.Range('D12').Value = 58
$ws.Range('A13').Value = '03:52:04'
$ws.Range('B13').Value = '05:16'
$ws.Range('C13').Value = '17_ROMERO'
$ws.Range('D13').Value = 84
$ws.Range('A14').Value = '04:32:18'
$ws.Range('B14').Value = '05:22'
$ws.Range('C14').Value = '23_HERNANDEZ'
$ws.Range('D14').Value = 50
$ws.Range('A15').Value = '04:13:31'
$ws.Range('B15').Value = '05:31'
$ws.Range('C15').Value = '81_EL PELIGRO'
$ws.Range('D15').Value = 78
$ws.Range('A16').Value = '04:32:18'
$ws.Range('B16').Value = '05:32'
$ws.Range('C16').Value = '81_EL PELIGRO'
$ws.Range('D16').Value = 60
$ws.Range('A17').Value = '03:52:04'
$ws.Range('B17').Value = '05:35'
$ws.Range('C17').Value = '215B_EL PATO'
$ws.Range('D17').Value = 103
$ws.Range('A18').Value = '03:52:04'
$ws.Range('B18').Value = '05:46'
$ws.Range('C18').Value = '15_ABASTO'
$ws.Range('D18').Value = 114
$ws.Range('A19').Value = '04:32:18'
$ws.Range('B19').Value = '05:47'
$ws.Range('C19').Value = '14_ABASTO'
$ws.Range('D19').Value = 75
$ws.Range('B20').Value = '05:50'
$ws.Range('C20').Value = '14_ABASTO'
$ws.Range('D20').Value = 97
$ws.Range('A21').Value = '04:32:18'
$ws.Range('B21').Value = '05:52'
$ws.Range('C21').Value = '17_ROMERO'
$ws.Range('D21').Value = 80
$ws.Range('E21').Value = 'LP1912'
$ws.Range('A22').Value = '04:32:18'
$ws.Range('B22').Value = '06:01'
$ws.Range('C22').Value = '16_SANTA ANA'
$ws.Range('D22').Value = 89
$ws.Range('E22').Value = 'LP1912'
$ws.Range('A23').Value = '04:13:31'
$ws.Range('B23').Value = '06:03'
$ws.Range('C23').Value = '10_OLMOS'
$ws.Range('D23').Value = 110
$ws.Range('E23').Value = 'LP1912'
$ws.Range('A24').Value = '04:32:18'
$ws.Range('B24').Value = '06:04'
$ws.Range('C24').Value = '10_OLMOS'
$ws.Range('D24').Value = 92
$ws.Range('E24').Value = 'LP1912'
$ws.Range('A25').Value = '04:32:18'
$ws.Range('B25').Value = '06:11'
$ws.Range('C25').Value = '215A_EL PATO'
$ws.Range('D25').Value = 99
$ws.Range('E25').Value = 'LP1912'
$ws.Range('A26').Value = '04:32:18'
$ws.Range('B26').Value = '06:15'
$ws.Range('C26').Value = '17_ROMERO'
$ws.Range('D26').Value = 103
$ws.Range('E26').Value = 'LP1912'
$ws.Range('A27').Value = '04:32:18'
$ws.Range('B27').Value = '06:24'
$ws.Range('C27').Value = '11_ETCHEVERRY'
$ws.Range('D27').Value = 112
$ws.Range('E27').Value = 'LP1912'
$ws.Range('A28').Value = '04:32:18'
$ws.Range('B28').Value = '06:27'
$ws.Range('C28').Value = '23_HERNANDEZ'
$ws.Range('D28').Value = 115
$ws.Range('E28').Value = 'LP1912'
$ws.Range('A29').Value = '04:32:18'
$ws.Range('B29').Value = '06:31'
$ws.Range('C29').Value = '16_SANTA ANA'
$ws.Range('D29').Value = 119
$ws.Range('E29').Value = 'LP1912'
$ws.Range('A30').Value = '04:32:18'
$ws.Range('B30').Value = '06:31'
$ws.Range('C30').Value = '17X38_ROMERO'
$ws.Range('D30').Value = 119
$ws.Range('E30').Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 04:32:18'
$ws.Range('A3').Value = 'Total filas: 5'
$ws.Range('A8').Value = '04:32:18'
$ws.Range('B8').Value = '04:47'
$ws.Range('C8').Value = '215_EL PELIGRO'
$ws.Range('D8').Value = 15
$ws.Range('A9').Value = '03:52:04'
$ws.Range('B9').Value = '05:35'
$ws.Range('C9').Value = '215B_EL PATO'
$ws.Range('D9').Value = 103
$ws.Range('A10').Value = '04:32:18'
$ws.Range('B10').Value = '06:11'
$ws.Range('C10').Value = '215A_EL PATO'
$ws.Range('D10').Value = 99
$ws.Range('E10').Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 04:32:18'
